$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New "NOTES" block (rows 9-14) describing the dataset columns.
# ---------------------------------------------------------------------------

$ws.Range("C9").Value = "NOTES"

$ws.Range("C10").Value = "d_per_hari"
$ws.Range("D10:F10").Merge()
$ws.Range("D10").Value = "Permintaan rata-rata per hari"

$ws.Range("C11").Value = "D_per_hari"
$ws.Range("D11:F11").Merge()
$ws.Range("D11").Value = "Permintaan total per tahun (biasanya = d_per_hari × 365)"

$ws.Range("C12").Value = "S"
$ws.Range("D12:F12").Merge()
$ws.Range("D12").Value = "Biaya pemesanan setiap kali order (ordering cost per order)"

$ws.Range("C13").Value = "h(Rp/unit/year)"
$ws.Range("D13:F13").Merge()
$ws.Range("D13").Value = "Biaya penyimpanan per unit per tahun (holding cost per unit per year)"

$ws.Range("C14").Value = "LeadTime(days)"
$ws.Range("D14:F14").Merge()
$ws.Range("D14").Value = "Waktu tunggu dari pemesanan hingga barang datang"

# Touch every merged D:F block in one go so each underlying cell (D,E,F) of
# every notes row is written out explicitly with a shared, plain style -
# mirroring the source workbook where the whole notes block uses one style.
$notes = $excel.Union($ws.Range("D10:F10"), $ws.Range("D11:F11"), $ws.Range("D12:F12"), $ws.Range("D13:F13"), $ws.Range("D14:F14"))
$notes.WrapText = $false

# ---------------------------------------------------------------------------
# Header row A1:F1 gets an accent-6 themed fill.
# ---------------------------------------------------------------------------

$ws.Range("A1:F1").Interior.ThemeColor = 10

# ---------------------------------------------------------------------------
# Column widths: A and B are now explicitly sized.
# ---------------------------------------------------------------------------

$ws.Columns.Item(1).ColumnWidth = 12.666666666666666
$ws.Columns.Item(2).ColumnWidth = 9.833333333333332

# ---------------------------------------------------------------------------
# Final selection, matching the saved workbook view.
# ---------------------------------------------------------------------------

$ws.Range("D14:F14").Select()
